# add results from latest run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast/revision figures for existing rows (2-11) with the latest model run
$ws.Cells.Item(2, 2).Value = 0.29274225186733338
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0

$ws.Cells.Item(3, 2).Value = 0.29250758799098536
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = -0.0035140448593851004
$ws.Cells.Item(3, 5).Value = 0.000035326584093716963
$ws.Cells.Item(3, 6).Value = -0.00066015558207592614
$ws.Cells.Item(3, 7).Value = 0.00049988994216002231
$ws.Cells.Item(3, 8).Value = -0.000036386518782310113
$ws.Cells.Item(3, 9).Value = -0.0010923731076776248
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.0000032390195756626383

$ws.Cells.Item(4, 2).Value = 0.29156153844083182
$ws.Cells.Item(4, 3).Value = -0.001382812978494654
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.000066371767707815646
$ws.Cells.Item(4, 6).Value = 0.000025190818588198996
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0.000096140559908248845
$ws.Cells.Item(4, 9).Value = -0.001402460898146023
$ws.Cells.Item(4, 10).Value = 0.000059092953741710421
$ws.Cells.Item(4, 11).Value = -0.000060482474875112047

$ws.Cells.Item(5, 2).Value = 0.29426080244315872
$ws.Cells.Item(5, 3).Value = 0.0045630010008605908
$ws.Cells.Item(5, 4).Value = -0.0053301109160312499
$ws.Cells.Item(5, 5).Value = 0.00000065397367624651922
$ws.Cells.Item(5, 6).Value = -0.00052057893485502515
$ws.Cells.Item(5, 7).Value = -0.0011817305155684176
$ws.Cells.Item(5, 8).Value = -0.000052123884068082571
$ws.Cells.Item(5, 9).Value = -0.00036248549665453392
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = -0.00015588972845836935

$ws.Cells.Item(6, 2).Value = 0.31905150050067566
$ws.Cells.Item(6, 3).Value = 0.026583263228862328
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = -0.00016097790358777427
$ws.Cells.Item(6, 6).Value = -0.000044146192946365815
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = -0.000071770639724629456
$ws.Cells.Item(6, 9).Value = -0.0022757391783104642
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = -0.0021262305811587945

$ws.Cells.Item(7, 2).Value = 0.31754139480399807
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = -0.001813426596032346
$ws.Cells.Item(7, 5).Value = -0.00032355386139361363
$ws.Cells.Item(7, 6).Value = -0.0014746742468685272
$ws.Cells.Item(7, 7).Value = 0.00078453919528978576
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.00017096731347514762
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = -0.0011796690648845365

$ws.Cells.Item(8, 2).Value = 0.27189410750104442
$ws.Cells.Item(8, 3).Value = -0.040302096336972446
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.000054282621489518804
$ws.Cells.Item(8, 6).Value = -0.00045333674438916547
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0.000037756818756590117
$ws.Cells.Item(8, 9).Value = 0.0016770117950774036
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = -0.00011208916280497361

$ws.Cells.Item(9, 2).Value = 0.26407735542389971
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0.0006047233673677206
$ws.Cells.Item(9, 5).Value = -0.0024311888866725313
$ws.Cells.Item(9, 6).Value = -0.0060083298392463865
$ws.Cells.Item(9, 7).Value = 0.00077294215961939939
$ws.Cells.Item(9, 8).Value = -0.00016016261892130829
$ws.Cells.Item(9, 9).Value = -0.00049121322006091871
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = -0.000019115905987576731

$ws.Cells.Item(10, 2).Value = 0.35321937438698242
$ws.Cells.Item(10, 3).Value = 0.073928888108657997
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = -0.000014657906929041119
$ws.Cells.Item(10, 6).Value = -0.00036043892303250626
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = -0.0000059425292918985265
$ws.Cells.Item(10, 9).Value = -0.00052027981559568215
$ws.Cells.Item(10, 10).Value = -0.0022613311529176981
$ws.Cells.Item(10, 11).Value = 0.0021471067489824414

$ws.Cells.Item(11, 2).Value = 0.34550359785410567
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = -0.011734236036316769
$ws.Cells.Item(11, 5).Value = 0.0013217972532459344
$ws.Cells.Item(11, 6).Value = 0.0014838154244163674
$ws.Cells.Item(11, 7).Value = 0.0013375142122967398
$ws.Cells.Item(11, 8).Value = 0.00024630087506194686
$ws.Cells.Item(11, 9).Value = -0.0042937141760055718
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0.0045527243814054708

# Append new row 12 for the 2025-08-30 survey vintage
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2025-08-30"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).Value = 0.16304185609794236
$ws.Cells.Item(12, 3).Value = -0.048246572967269252
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0.000099733070870745703
$ws.Cells.Item(12, 6).Value = 0.0000068780163400907228
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0.00000068315401009552564
$ws.Cells.Item(12, 9).Value = -0.00028509358276941854
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = -0.0075302770571035182
